$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 88

# Force text entry (matching the existing column's literal-string style) by
# prefixing with an apostrophe, so Excel doesn't auto-coerce the date/time
# looking strings into date/time serial values.
$ws.Cells.Item($newRow, 1).Value = "'2025-10-19"
$ws.Cells.Item($newRow, 2).Value = "'21:20:13"
$ws.Cells.Item($newRow, 3).Value = "1.00 EUR = 1,703.9130"
